# Applies the "testing of bot done." commit:
#  - refreshes week-04 -> week-07 delivery paths used on the Settings sheet
#  - adds two new Settings keys (CorrectnameFilePath / CorrectSheetName)
#  - adds a large block of new mail-notification / working-date constants
#    on the Constants sheet
#  - leaves the workbook with the Constants sheet as the active tab,
#    matching where the author was working when they saved

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Settings sheet: refresh the week 04 -> week 07 file-share paths
# ---------------------------------------------------------------------
$pathBaseDatos   = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/GULF/Base de Datos"
$pathVip         = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/GULF/VIP"
$pathTradicional = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/GULF/TRADICIONAL"
$pathExportadas  = "/Planeacion/0.Envios TS/2022/07 Envío Semana 07/GULF/Base de Datos/Exportadas"

$ws1.Range("B2").Value  = $pathBaseDatos     # InputnonPartnerFilePath
$ws1.Range("B7").Value  = $pathVip           # BDVIPfilePath
$ws1.Range("B9").Value  = $pathTradicional   # BDTempFilePath
$ws1.Range("B12").Value = $pathBaseDatos     # BlackListTrackingPath
$ws1.Range("B15").Value = $pathExportadas    # SurveyTradFilePath
$ws1.Range("B18").Value = $pathExportadas    # SurveyVipFilePath

# Two new settings appended at the bottom of the sheet
$ws1.Range("A27").Value = "CorrectnameFilePath"
$ws1.Range("A28").Value = "CorrectSheetName"

$ws1.Columns("B").AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Constants sheet: append the new notification / working-date constants
# ---------------------------------------------------------------------

# Working week boundaries
$ws2.Range("A6").Value = "DataNoBelongToCurrentWeekMailSubject"
$ws2.Range("B6").Value = "Notification- Data doesn't belong to this week"

$ws2.Range("A7").Value = "DataNoBelongToCurrentWeekMailBody"
$dataNoBelongBody = @"
Dear Team<br/>
Please find the below data which does not belong to current week.<br/>
[Nothisweekdatatable]<br>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws2.Range("B7").Value = $dataNoBelongBody
$ws2.Range("B7").WrapText = $true

$ws2.Range("A9").Value = "MailTo"
$ws2.Range("B9").Value = "rhernandez@tecnoyar.com.mx"

$ws2.Range("A10").Value = "MailCC"
$ws2.Range("B10").Value = "arerodriguez@tecnoyar.com"

$ws2.Range("A12").Value = "InputFileExceptionSubject"
$ws2.Range("B12").Value = "Error- Input file Notification"

$ws2.Range("A13").Value = "InputFileExceptionbody"
$inputFileExceptionBody = @"
Dear Team<br/>
There is a error in Cinepolish Gulf run, below is the detail of error<br/>
[error]<br/>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws2.Range("B13").Value = $inputFileExceptionBody
$ws2.Range("B13").WrapText = $true

$ws2.Range("A15").Value = "downloadFileExceptionSubject"
$ws2.Range("B15").Value = "Notification- Not All file Got downloaded"

$ws2.Range("A16").Value = "downloadFileExceptionbody"
$downloadFileExceptionBody = @"
Dear Team<br/>
There is a error in Cinepos Gulf process run, below is the detail of error<br/>
Not all required files got downloaded from FTP<br/>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws2.Range("B16").Value = $downloadFileExceptionBody
$ws2.Range("B16").WrapText = $true

$ws2.Range("A18").Value = "WrongNameMailSubject"
$ws2.Range("B18").Value = "Notification- name is not correct"

$ws2.Range("A19").Value = "WrongNameMailbody"
$wrongNameMailBody = @"
Dear Team<br/>
below names are not correct in Gulf process input file.<br/>
[wrongnamedatatable]<br>
Also let us know if anything is required<br/>
Thank you,<br/>
**********************This is system generated E-Mail, please do not respond on this************
"@
$ws2.Range("B19").Value = $wrongNameMailBody
$ws2.Range("B19").WrapText = $true

$ws2.Range("A21").Value = "NADirectrioMailSubject"
$ws2.Range("B21").Value = "file of the Maxico mark contains errors in the ""N"" column called ""clave_tipo_cine"""

$ws2.Range("A22").Value = "NADirectrioMailBody"
$ws2.Range("B22").Value = "#NA data found at input file in ""N"" column"

$ws2.Range("A24").Value = "WorkingStartDate"
$ws2.Range("B24").Value = (Get-Date -Year 2022 -Month 2 -Day 7).Date
$ws2.Range("B24").NumberFormat = "mm-dd-yy"
$ws2.Range("C24").Value = "dd/MM/yyyy"
$ws2.Range("C24").NumberFormat = "mm-dd-yy"

$ws2.Range("A25").Value = "WorkingEndDate"
$ws2.Range("B25").Value = (Get-Date -Year 2022 -Month 2 -Day 13).Date
$ws2.Range("B25").NumberFormat = "mm-dd-yy"
$ws2.Range("C25").Value = "dd/MM/yyyy"
$ws2.Range("C25").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------
# Leave the workbook focused on the Constants sheet, near where the new
# rows were added, matching the saved view state.
# ---------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("B94").Select() | Out-Null
